$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (datetime serial, symbol, open, high, low, close, volume)
$newRows = @(
    @(44986.45833333334, "ECONOMICS:NGM2", 54191661320000, 54191661320000, 54191661320000, 54191661320000, 0),
    @(45017.45833333334, "ECONOMICS:NGM2", 55646745840000, 55646745840000, 55646745840000, 55646745840000, 0),
    @(45047.41666666666, "ECONOMICS:NGM2", 55500913410000, 55500913410000, 55500913410000, 55500913410000, 0)
)

$lastRow = 279
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $lastRow + 1 + $i
    $srcRange = $ws.Range("A" + $lastRow + ":G" + $lastRow)
    $dstRange = $ws.Range("A" + $targetRow + ":G" + $targetRow)

    # Copy formatting (style/number format/border/alignment) from the last data row
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)

    $rowData = $newRows[$i]
    $ws.Cells.Item($targetRow, 1).Value = $rowData[0]
    $ws.Cells.Item($targetRow, 2).Value = $rowData[1]
    $ws.Cells.Item($targetRow, 3).Value = $rowData[2]
    $ws.Cells.Item($targetRow, 4).Value = $rowData[3]
    $ws.Cells.Item($targetRow, 5).Value = $rowData[4]
    $ws.Cells.Item($targetRow, 6).Value = $rowData[5]
    $ws.Cells.Item($targetRow, 7).Value = $rowData[6]
}

$wb.Save()
